# Helper: force a cell to be stored as TEXT (matches source data where
# fund codes / percentages / amounts are kept as text, not numbers).
function Set-TextCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Helper: set a cell as a genuine number.
function Set-NumCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a brand new "2022-Q4" sheet right after "总计", by cloning
#    the structure/formatting of the existing "2022-Q3" sheet (same
#    column layout & header style) and then overwriting its content.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$template = $wb.Worksheets("2022-Q3")

$template.Copy([System.Reflection.Missing]::Value, $total)
$newQ4 = $wb.Worksheets.Item(2)
$newQ4.Name = "2022-Q4"

# The template (old "2022-Q3") has 26 data rows (rows 2-27); the new
# "2022-Q4" sheet only needs 23 data rows (rows 2-24). Remove the extra
# trailing rows 25-27.
$newQ4.Range("A25:H27").EntireRow.Delete()

# ---------------------------------------------------------------------
# 2) Fill in the real 2022-Q4 fund-holding data (overwrites whatever
#    was copied over from the "2022-Q3" template).
# ---------------------------------------------------------------------
Set-TextCell $newQ4 1 2 '基金代码'
Set-TextCell $newQ4 1 3 '基金名称'
Set-TextCell $newQ4 1 4 '基金规模'
Set-TextCell $newQ4 1 5 '股票总仓位'
Set-TextCell $newQ4 1 6 '仓位占比'
Set-TextCell $newQ4 1 7 '持有市值(亿元)'
Set-TextCell $newQ4 1 8 '仓位排名'
Set-NumCell $newQ4 2 1 0
Set-TextCell $newQ4 2 2 '003231'
Set-TextCell $newQ4 2 3 '创金合信医疗保健行业股票C'
Set-TextCell $newQ4 2 4 '6.48'
Set-TextCell $newQ4 2 5 '94.76'
Set-TextCell $newQ4 2 6 '8.43'
Set-TextCell $newQ4 2 7 '0.5463'
Set-NumCell $newQ4 2 8 3
Set-NumCell $newQ4 3 1 1
Set-TextCell $newQ4 3 2 '003230'
Set-TextCell $newQ4 3 3 '创金合信医疗保健行业股票A'
Set-TextCell $newQ4 3 4 '4.22'
Set-TextCell $newQ4 3 5 '94.76'
Set-TextCell $newQ4 3 6 '8.43'
Set-TextCell $newQ4 3 7 '0.3557'
Set-NumCell $newQ4 3 8 3
Set-NumCell $newQ4 4 1 2
Set-TextCell $newQ4 4 2 '009960'
Set-TextCell $newQ4 4 3 '银华多元机遇混合'
Set-TextCell $newQ4 4 4 '9.55'
Set-TextCell $newQ4 4 5 '90.95'
Set-TextCell $newQ4 4 6 '3.42'
Set-TextCell $newQ4 4 7 '0.3266'
Set-NumCell $newQ4 4 8 2
Set-NumCell $newQ4 5 1 3
Set-TextCell $newQ4 5 2 '010585'
Set-TextCell $newQ4 5 3 '创金合信医药消费股票A'
Set-TextCell $newQ4 5 4 '3.60'
Set-TextCell $newQ4 5 5 '93.18'
Set-TextCell $newQ4 5 6 '7.35'
Set-TextCell $newQ4 5 7 '0.2646'
Set-NumCell $newQ4 5 8 5
Set-NumCell $newQ4 6 1 4
Set-TextCell $newQ4 6 2 '013067'
Set-TextCell $newQ4 6 3 '富安达中小盘六个月持有期混合'
Set-TextCell $newQ4 6 4 '2.09'
Set-TextCell $newQ4 6 5 '78.87'
Set-TextCell $newQ4 6 6 '7.45'
Set-TextCell $newQ4 6 7 '0.1557'
Set-NumCell $newQ4 6 8 1
Set-NumCell $newQ4 7 1 5
Set-TextCell $newQ4 7 2 '320012'
Set-TextCell $newQ4 7 3 '诺安主题精选混合'
Set-TextCell $newQ4 7 4 '3.57'
Set-TextCell $newQ4 7 5 '90.77'
Set-TextCell $newQ4 7 6 '3.97'
Set-TextCell $newQ4 7 7 '0.1417'
Set-NumCell $newQ4 7 8 10
Set-NumCell $newQ4 8 1 6
Set-TextCell $newQ4 8 2 '014737'
Set-TextCell $newQ4 8 3 '创金合信专精特新股票C'
Set-TextCell $newQ4 8 4 '1.83'
Set-TextCell $newQ4 8 5 '91.82'
Set-TextCell $newQ4 8 6 '7.04'
Set-TextCell $newQ4 8 7 '0.1288'
Set-NumCell $newQ4 8 8 6
Set-NumCell $newQ4 9 1 7
Set-TextCell $newQ4 9 2 '013349'
Set-TextCell $newQ4 9 3 '创金合信大健康混合C'
Set-TextCell $newQ4 9 4 '1.38'
Set-TextCell $newQ4 9 5 '94.19'
Set-TextCell $newQ4 9 6 '9.04'
Set-TextCell $newQ4 9 7 '0.1248'
Set-NumCell $newQ4 9 8 2
Set-NumCell $newQ4 10 1 8
Set-TextCell $newQ4 10 2 '310388'
Set-TextCell $newQ4 10 3 '申万菱信消费增长混合A'
Set-TextCell $newQ4 10 4 '3.09'
Set-TextCell $newQ4 10 5 '92.93'
Set-TextCell $newQ4 10 6 '3.39'
Set-TextCell $newQ4 10 7 '0.1048'
Set-NumCell $newQ4 10 8 9
Set-NumCell $newQ4 11 1 9
Set-TextCell $newQ4 11 2 '011383'
Set-TextCell $newQ4 11 3 '富安达医药创新混合'
Set-TextCell $newQ4 11 4 '1.58'
Set-TextCell $newQ4 11 5 '89.59'
Set-TextCell $newQ4 11 6 '6.61'
Set-TextCell $newQ4 11 7 '0.1044'
Set-NumCell $newQ4 11 8 1
Set-NumCell $newQ4 12 1 10
Set-TextCell $newQ4 12 2 '013348'
Set-TextCell $newQ4 12 3 '创金合信大健康混合A'
Set-TextCell $newQ4 12 4 '0.82'
Set-TextCell $newQ4 12 5 '94.19'
Set-TextCell $newQ4 12 6 '9.04'
Set-TextCell $newQ4 12 7 '0.0741'
Set-NumCell $newQ4 12 8 2
Set-NumCell $newQ4 13 1 11
Set-TextCell $newQ4 13 2 '014736'
Set-TextCell $newQ4 13 3 '创金合信专精特新股票A'
Set-TextCell $newQ4 13 4 '0.81'
Set-TextCell $newQ4 13 5 '91.82'
Set-TextCell $newQ4 13 6 '7.04'
Set-TextCell $newQ4 13 7 '0.0570'
Set-NumCell $newQ4 13 8 6
Set-NumCell $newQ4 14 1 12
Set-TextCell $newQ4 14 2 '005108'
Set-TextCell $newQ4 14 3 '圆信永丰双利优选定期开放灵活配置混合'
Set-TextCell $newQ4 14 4 '0.73'
Set-TextCell $newQ4 14 5 '88.26'
Set-TextCell $newQ4 14 6 '5.26'
Set-TextCell $newQ4 14 7 '0.0384'
Set-NumCell $newQ4 14 8 5
Set-NumCell $newQ4 15 1 13
Set-TextCell $newQ4 15 2 '001861'
Set-TextCell $newQ4 15 3 '富安达健康人生灵活配置混合A'
Set-TextCell $newQ4 15 4 '0.55'
Set-TextCell $newQ4 15 5 '91.87'
Set-TextCell $newQ4 15 6 '6.98'
Set-TextCell $newQ4 15 7 '0.0384'
Set-NumCell $newQ4 15 8 1
Set-NumCell $newQ4 16 1 14
Set-TextCell $newQ4 16 2 '010586'
Set-TextCell $newQ4 16 3 '创金合信医药消费股票C'
Set-TextCell $newQ4 16 4 '0.52'
Set-TextCell $newQ4 16 5 '93.18'
Set-TextCell $newQ4 16 6 '7.35'
Set-TextCell $newQ4 16 7 '0.0382'
Set-NumCell $newQ4 16 8 5
Set-NumCell $newQ4 17 1 15
Set-TextCell $newQ4 17 2 '001965'
Set-TextCell $newQ4 17 3 '圆信永丰兴源灵活配置混合A'
Set-TextCell $newQ4 17 4 '0.58'
Set-TextCell $newQ4 17 5 '84.58'
Set-TextCell $newQ4 17 6 '5.27'
Set-TextCell $newQ4 17 7 '0.0306'
Set-NumCell $newQ4 17 8 4
Set-NumCell $newQ4 18 1 16
Set-TextCell $newQ4 18 2 '015254'
Set-TextCell $newQ4 18 3 '申万菱信消费增长混合C'
Set-TextCell $newQ4 18 4 '0.89'
Set-TextCell $newQ4 18 5 '92.93'
Set-TextCell $newQ4 18 6 '3.39'
Set-TextCell $newQ4 18 7 '0.0302'
Set-NumCell $newQ4 18 8 9
Set-NumCell $newQ4 19 1 17
Set-TextCell $newQ4 19 2 '006274'
Set-TextCell $newQ4 19 3 '圆信永丰医药健康混合'
Set-TextCell $newQ4 19 4 '0.24'
Set-TextCell $newQ4 19 5 '85.59'
Set-TextCell $newQ4 19 6 '5.43'
Set-TextCell $newQ4 19 7 '0.0130'
Set-NumCell $newQ4 19 8 4
Set-NumCell $newQ4 20 1 18
Set-TextCell $newQ4 20 2 '012019'
Set-TextCell $newQ4 20 3 '国投瑞银安泽混合A'
Set-TextCell $newQ4 20 4 '0.62'
Set-TextCell $newQ4 20 5 '31.81'
Set-TextCell $newQ4 20 6 '0.94'
Set-TextCell $newQ4 20 7 '0.0058'
Set-NumCell $newQ4 20 8 10
Set-NumCell $newQ4 21 1 19
Set-TextCell $newQ4 21 2 '001966'
Set-TextCell $newQ4 21 3 '圆信永丰兴源灵活配置混合C'
Set-TextCell $newQ4 21 4 '0.09'
Set-TextCell $newQ4 21 5 '84.58'
Set-TextCell $newQ4 21 6 '5.27'
Set-TextCell $newQ4 21 7 '0.0047'
Set-NumCell $newQ4 21 8 4
Set-NumCell $newQ4 22 1 20
Set-TextCell $newQ4 22 2 '562520'
Set-TextCell $newQ4 22 3 '华夏中证智选1000成长创新策略ETF'
Set-TextCell $newQ4 22 4 '0.38'
Set-TextCell $newQ4 22 5 '96.24'
Set-TextCell $newQ4 22 6 '1.05'
Set-TextCell $newQ4 22 7 '0.0040'
Set-NumCell $newQ4 22 8 3
Set-NumCell $newQ4 23 1 21
Set-TextCell $newQ4 23 2 '012020'
Set-TextCell $newQ4 23 3 '国投瑞银安泽混合C'
Set-TextCell $newQ4 23 4 '0.11'
Set-TextCell $newQ4 23 5 '31.81'
Set-TextCell $newQ4 23 6 '0.94'
Set-TextCell $newQ4 23 7 '0.0010'
Set-NumCell $newQ4 23 8 10
Set-NumCell $newQ4 24 1 22
Set-TextCell $newQ4 24 2 '014470'
Set-TextCell $newQ4 24 3 '富安达健康人生灵活配置混合C'
Set-TextCell $newQ4 24 4 '0.00'
Set-TextCell $newQ4 24 5 '91.87'
Set-TextCell $newQ4 24 6 '6.98'
Set-NumCell $newQ4 24 7 0
Set-NumCell $newQ4 24 8 1

# ---------------------------------------------------------------------
# 3) Update the "总计" (Total) summary sheet: insert a new row for
#    2022-Q4 right after the header row, push the existing rows down,
#    and fix up the running index in column A.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

Set-NumCell  $total 2 1 0
Set-TextCell $total 2 2 '2022-Q4'
Set-NumCell  $total 2 3 23
Set-NumCell  $total 2 4 2.59

# Re-number column A (running index) for the rows that shifted down by
# one position (old rows 2-9 are now rows 3-10, and need index 1-8
# instead of the 0-7 they kept after the shift).
for ($r = 3; $r -le 10; $r++) {
    $idx = $r - 2
    Set-NumCell $total $r 1 $idx
}

Write-Host "Done"
